# Updated cryptos list values (price/volume refresh + a few row reorderings)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.982.45"
$ws.Range("E2").Value = "  +0.58%  "

$ws.Range("D3").Value = "2.301.45"
$ws.Range("E3").Value = "  +0.52%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "114.36"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +18.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "270.90"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.46%  "

$ws.Range("E7").Value = "  +1.81%  "

$ws.Range("E8").Value = "  +0.29%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.623"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.53%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.23"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +6.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0950"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.60%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.09"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +14.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.107"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.90"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.36%  "

$ws.Range("D15").Value = "2.647.11"
$ws.Range("E15").Value = "  +0.53%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.855"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.15%  "

$ws.Range("D17").Value = "2.301.53"
$ws.Range("E17").Value = "  +0.57%  "

$ws.Range("D18").Value = "43.859.15"
$ws.Range("E18").Value = "  +0.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000111"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.55%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.82"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +10.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.61"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.46"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.39"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.32%  "

$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.95"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +8.85%  "

$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.69"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +6.34%  "

$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.74"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.97%  "

$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "41.94"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +7.96%  "

$ws.Range("E29").Value = "  -1.91%  "

$ws.Range("E30").Value = "  -0.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.59"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.47%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0938"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.65"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.72"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +5.42%  "

$ws.Range("E35").Value = "  +1.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.67"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.74%  "

$ws.Range("E37").Value = "  +4.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.108"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.62%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.82"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +7.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "74.36"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +15.54%  "

$ws.Range("E41").Value = "  +3.62%  "

$ws.Range("B42").Value = "THORChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.44"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +23.81%  "

$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.71"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +10.75%  "

$ws.Range("E44").Value = "  +3.37%  "

$ws.Range("E45").Value = "  +0.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.40"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.13%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.88"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.84%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.77"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.50%  "

$ws.Range("E49").Value = "  -2.00%  "

$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.472"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +9.37%  "

$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.24"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.80%  "
